$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.34"
$ws.Range("E2").Value = "'4.94%"
$ws.Range("G2").Value = "'17"
$ws.Range("D3").Value = "'28.00"
$ws.Range("E3").Value = "'-3.87%"
$ws.Range("G3").Value = "'17"
$ws.Range("E4").Value = "'-1.05%"
$ws.Range("G4").Value = "'17"
$ws.Range("D5").Value = "'0.05902"
$ws.Range("E5").Value = "'3.38%"
$ws.Range("G5").Value = "'17"
$ws.Range("D6").Value = "'6.708"
$ws.Range("E6").Value = "'1.32%"
$ws.Range("G6").Value = "'17"
$ws.Range("D7").Value = "'0.8711"
$ws.Range("E7").Value = "'2.28%"
$ws.Range("G7").Value = "'17"
$ws.Range("D8").Value = "'0.9995"
$ws.Range("E8").Value = "'16.60%"
$ws.Range("G8").Value = "'17"
$ws.Range("D9").Value = "'0.1411"
$ws.Range("E9").Value = "'2.91%"
$ws.Range("G9").Value = "'17"
$ws.Range("D10").Value = "'0.07203"
$ws.Range("E10").Value = "'2.26%"
$ws.Range("G10").Value = "'17"
$ws.Range("D11").Value = "'0.03159"
$ws.Range("E11").Value = "'-0.07%"
$ws.Range("G11").Value = "'17"
$ws.Range("D12").Value = "'0.09219"
$ws.Range("E12").Value = "'-0.69%"
$ws.Range("G12").Value = "'17"
$ws.Range("D13").Value = "'0.001551"
$ws.Range("E13").Value = "'1.82%"
$ws.Range("G13").Value = "'17"
$ws.Range("D14").Value = "'0.0006078"
$ws.Range("E14").Value = "'1.70%"
$ws.Range("G14").Value = "'17"
$ws.Range("D15").Value = "'0.005855"
$ws.Range("E15").Value = "'-4.00%"
$ws.Range("G15").Value = "'17"
$ws.Range("D16").Value = "'3.494"
$ws.Range("E16").Value = "'0.05%"
$ws.Range("G16").Value = "'17"
$ws.Range("E17").Value = "'1.50%"
$ws.Range("G17").Value = "'17"
$ws.Range("E18").Value = "'0.06%"
$ws.Range("G18").Value = "'17"
$ws.Range("D19").Value = "'0.3123"
$ws.Range("E19").Value = "'-1.12%"
$ws.Range("G19").Value = "'17"
$ws.Range("D20").Value = "'0.03644"
$ws.Range("E20").Value = "'9.69%"
$ws.Range("G20").Value = "'17"
$ws.Range("E21").Value = "'0.96%"
$ws.Range("G21").Value = "'17"
$ws.Range("D22").Value = "'3.525"
$ws.Range("E22").Value = "'1.07%"
$ws.Range("G22").Value = "'17"
$ws.Range("D23").Value = "'0.04201"
$ws.Range("E23").Value = "'2.51%"
$ws.Range("G23").Value = "'17"
$ws.Range("D24").Value = "'0.1362"
$ws.Range("E24").Value = "'-1.25%"
$ws.Range("G24").Value = "'17"
$ws.Range("D25").Value = "'0.001217"
$ws.Range("E25").Value = "'-0.40%"
$ws.Range("G25").Value = "'17"
$ws.Range("E26").Value = "'9.81%"
$ws.Range("G26").Value = "'17"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("E27").Value = "'-0.03%"
$ws.Range("G27").Value = "'17"
$ws.Range("D28").Value = "'0.0001937"
$ws.Range("E28").Value = "'33.56%"
$ws.Range("G28").Value = "'17"
$ws.Range("G29").Value = "'17"
$ws.Range("G30").Value = "'17"
$ws.Range("G31").Value = "'17"
$ws.Range("G32").Value = "'17"
$ws.Range("G33").Value = "'17"
$ws.Range("G34").Value = "'17"
$ws.Range("G35").Value = "'17"
$ws.Range("G36").Value = "'17"
$ws.Range("G37").Value = "'17"
$ws.Range("G38").Value = "'17"
$ws.Range("G39").Value = "'17"
$ws.Range("D40").Value = "'0.03837"
$ws.Range("E40").Value = "'2.17%"
$ws.Range("G40").Value = "'17"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1105"
$ws.Range("E41").Value = "'3.73%"
$ws.Range("G41").Value = "'17"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003994"
$ws.Range("E42").Value = "'-22.17%"
$ws.Range("G42").Value = "'17"
$ws.Range("D43").Value = "'0.002449"
$ws.Range("E43").Value = "'-0.02%"
$ws.Range("G43").Value = "'17"
$ws.Range("E44").Value = "'14.00%"
$ws.Range("G44").Value = "'17"
$ws.Range("D45").Value = "'0.00005422"
$ws.Range("E45").Value = "'2.53%"
$ws.Range("G45").Value = "'17"
$ws.Range("E46").Value = "'-0.15%"
$ws.Range("G46").Value = "'17"
$ws.Range("D47").Value = "'0.08547"
$ws.Range("E47").Value = "'13.83%"
$ws.Range("G47").Value = "'17"
$ws.Range("D48").Value = "'0.002137"
$ws.Range("E48").Value = "'-12.60%"
$ws.Range("G48").Value = "'17"
$ws.Range("E49").Value = "'-0.15%"
$ws.Range("G49").Value = "'17"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("G50").Value = "'17"
$ws.Range("G51").Value = "'17"
